$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.11 = 7719.11 pesos`n✅ 7719.11 pesos = 2.11 = 930.01 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 472.995
$ws2.Range("O10").Value = 3651.1
$ws2.Range("N12").Value = 3652
$ws2.Range("O12").Value = 440
